$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Inflammatory-Mac"
$ws.Range("B2").Value = "Il12b"
$ws.Range("C2").Value = "Il12rb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.284858666666667
$ws.Range("H2").Value = 3.854576
$ws.Range("I2").Value = 0.6825120125588942
$ws.Range("J2").Value = 0.6825120125588942
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.2848286666666667
$ws.Range("N2").Value = 0.854486
$ws.Range("O2").Value = 0.08022967564521397
$ws.Range("P2").Value = 0.08022967564521397
$ws.Range("Q2").Value = 0.3659645808817779
$ws.Range("R2").Value = 3.293681227936
$ws.Range("S2").Value = 0.05475771739156228
$ws.Range("T2").Value = 0.05475771739156228
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("B3").Value = "Il12b"
$ws.Range("C3").Value = "Il12rb1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.284858666666667
$ws.Range("H3").Value = 3.854576
$ws.Range("I3").Value = 0.6825120125588942
$ws.Range("J3").Value = 0.6825120125588942
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.9694063333333333
$ws.Range("N3").Value = 2.908219
$ws.Range("O3").Value = 0.2730594381596053
$ws.Range("P3").Value = 0.2730594381596053
$ws.Range("Q3").Value = 1.245550128904889
$ws.Range("R3").Value = 11.209951160144
$ws.Range("S3").Value = 0.1863663466865131
$ws.Range("T3").Value = 0.1863663466865131
$ws.Range("A4").Value = "Inflammatory-Mac"
$ws.Range("B4").Value = "Il12b"
$ws.Range("C4").Value = "Il12rb1"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.284858666666667
$ws.Range("H4").Value = 3.854576
$ws.Range("I4").Value = 0.6825120125588942
$ws.Range("J4").Value = 0.6825120125588942
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.351483666666667
$ws.Range("N4").Value = 4.054451
$ws.Range("O4").Value = 0.3806818235166093
$ws.Range("P4").Value = 0.3806818235166093
$ws.Range("Q4").Value = 1.736465501975111
$ws.Range("R4").Value = 15.628189517776
$ws.Range("S4").Value = 0.2598199175129108
$ws.Range("T4").Value = 0.2598199175129108
$ws.Range("A5").Value = "Inflammatory-Mac"
$ws.Range("B5").Value = "Il12b"
$ws.Range("C5").Value = "Il12rb1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.284858666666667
$ws.Range("H5").Value = 3.854576
$ws.Range("I5").Value = 0.6825120125588942
$ws.Range("J5").Value = 0.6825120125588942
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.23571
$ws.Range("N5").Value = 0.70713
$ws.Range("O5").Value = 0.066394078474077
$ws.Range("P5").Value = 0.066394078474077
$ws.Range("Q5").Value = 0.30285403632
$ws.Range("R5").Value = 2.72568632688
$ws.Range("S5").Value = 0.04531475612133544
$ws.Range("T5").Value = 0.04531475612133544
$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("B6").Value = "Il12b"
$ws.Range("C6").Value = "Il12rb1"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 1.284858666666667
$ws.Range("H6").Value = 3.854576
$ws.Range("I6").Value = 0.6825120125588942
$ws.Range("J6").Value = 0.6825120125588942
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.7087373333333332
$ws.Range("N6").Value = 2.126212
$ws.Range("O6").Value = 0.1996349842044944
$ws.Range("P6").Value = 0.1996349842044945
$ws.Range("Q6").Value = 0.9106273051235555
$ws.Range("R6").Value = 8.195645746112
$ws.Range("S6").Value = 0.1362532748465725
$ws.Range("T6").Value = 0.1362532748465726
$ws.Range("A7").Value = "Resolving-Mac"
$ws.Range("B7").Value = "Il12b"
$ws.Range("C7").Value = "Il12rb1"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.597685
$ws.Range("H7").Value = 1.793055
$ws.Range("I7").Value = 0.3174879874411058
$ws.Range("J7").Value = 0.3174879874411058
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.2848286666666667
$ws.Range("N7").Value = 0.854486
$ws.Range("O7").Value = 0.08022967564521397
$ws.Range("P7").Value = 0.08022967564521397
$ws.Range("Q7").Value = 0.1702378216366667
$ws.Range("R7").Value = 1.53214039473
$ws.Range("S7").Value = 0.02547195825365169
$ws.Range("T7").Value = 0.02547195825365169
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Il12b"
$ws.Range("C8").Value = "Il12rb1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.597685
$ws.Range("H8").Value = 1.793055
$ws.Range("I8").Value = 0.3174879874411058
$ws.Range("J8").Value = 0.3174879874411058
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.9694063333333333
$ws.Range("N8").Value = 2.908219
$ws.Range("O8").Value = 0.2730594381596053
$ws.Range("P8").Value = 0.2730594381596053
$ws.Range("Q8").Value = 0.5793996243383334
$ws.Range("R8").Value = 5.214596619045
$ws.Range("S8").Value = 0.08669309147309219
$ws.Range("T8").Value = 0.08669309147309219
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Il12b"
$ws.Range("C9").Value = "Il12rb1"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.597685
$ws.Range("H9").Value = 1.793055
$ws.Range("I9").Value = 0.3174879874411058
$ws.Range("J9").Value = 0.3174879874411058
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.351483666666667
$ws.Range("N9").Value = 4.054451
$ws.Range("O9").Value = 0.3806818235166093
$ws.Range("P9").Value = 0.3806818235166093
$ws.Range("Q9").Value = 0.8077615153116667
$ws.Range("R9").Value = 7.269853637805
$ws.Range("S9").Value = 0.1208619060036985
$ws.Range("T9").Value = 0.1208619060036985
$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("B10").Value = "Il12b"
$ws.Range("C10").Value = "Il12rb1"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.597685
$ws.Range("H10").Value = 1.793055
$ws.Range("I10").Value = 0.3174879874411058
$ws.Range("J10").Value = 0.3174879874411058
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.23571
$ws.Range("N10").Value = 0.70713
$ws.Range("O10").Value = 0.066394078474077
$ws.Range("P10").Value = 0.066394078474077
$ws.Range("Q10").Value = 0.14088033135
$ws.Range("R10").Value = 1.26792298215
$ws.Range("S10").Value = 0.02107932235274155
$ws.Range("T10").Value = 0.02107932235274155
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Il12b"
$ws.Range("C11").Value = "Il12rb1"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.597685
$ws.Range("H11").Value = 1.793055
$ws.Range("I11").Value = 0.3174879874411058
$ws.Range("J11").Value = 0.3174879874411058
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.7087373333333332
$ws.Range("N11").Value = 2.126212
$ws.Range("O11").Value = 0.1996349842044944
$ws.Range("P11").Value = 0.1996349842044945
$ws.Range("Q11").Value = 0.4236016730733333
$ws.Range("R11").Value = 3.81241505766
$ws.Range("S11").Value = 0.06338170935792189
$ws.Range("T11").Value = 0.06338170935792191
